{"js": "// Update the reservation-form table: check-in time, room/cottage names\n// and all the fee figures now that the booking has gone through.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, column) -> new text, 0-based, matching the visible table grid.\nconst edits = [\n  [5, 1, \"2023-12-11 09:01\"], // CHECK IN DATE & TIME\n  [7, 1, \"Superior Room-1\"],  // ROOM name\n  [7, 2, \"2000\"],             // ROOM price\n  [9, 1, \"Kubo-3\"],           // COTTAGE name\n  [9, 2, \"1000\"],             // COTTAGE price\n  [10, 2, \"0.00\"],            // ADULTS price\n  [11, 2, \"0.00\"],            // KIDS price\n  [12, 2, \"0.00\"],            // SENIOR/PWD price\n  [13, 2, \"3000.00\"],         // TOTAL DUE\n  [14, 2, \"1500\"],            // DOWNPAYMENT\n  [16, 2, \"3000.00\"],         // TOTAL PAYABLE\n];\n\nfor (const [row, col, text] of edits) {\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the reservation-form table: check-in time, room/cottage names\n# and all the fee figures now that the booking has gone through.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(6, 2).Range.Text  = \"2023-12-11 09:01\"  # CHECK IN DATE & TIME\n$t.Cell(8, 2).Range.Text  = \"Superior Room-1\"   # ROOM name\n$t.Cell(8, 3).Range.Text  = \"2000\"              # ROOM price\n$t.Cell(10, 2).Range.Text = \"Kubo-3\"            # COTTAGE name\n$t.Cell(10, 3).Range.Text = \"1000\"              # COTTAGE price\n$t.Cell(11, 3).Range.Text = \"0.00\"              # ADULTS price\n$t.Cell(12, 3).Range.Text = \"0.00\"              # KIDS price\n$t.Cell(13, 3).Range.Text = \"0.00\"              # SENIOR/PWD price\n$t.Cell(14, 3).Range.Text = \"3000.00\"           # TOTAL DUE\n$t.Cell(15, 3).Range.Text = \"1500\"              # DOWNPAYMENT\n$t.Cell(17, 3).Range.Text = \"3000.00\"           # TOTAL PAYABLE\n"}
